# Updated cryptos list on Thu Jun 27 17:50:03 UTC 2024 with GitHub Actions
# Refreshes the Price (column D) and Volume(1h) (column E) figures for every
# coin row, and swaps three pairs of adjacent rows whose rank/ordering
# changed (Binance-PegBSC-USD <-> Fetch.AI, Stacks <-> FirstDigitalUSD,
# Maker <-> ONDO), updating their Coin/Link/Price/Volume values accordingly.
#
# Note: several new Price values are numeric-looking (e.g. "580.29"), but
# the source data stores Price as text (the original thousand separators
# use '.' the same as the decimal separator, e.g. "61.816.70"), so a
# leading "'" is used to force those particular cells to stay text instead
# of being auto-converted to a number by Excel.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = '61.780.18'
$ws.Cells.Item(2, 5).Value = '  +1.20%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '3.458.08'
$ws.Cells.Item(3, 5).Value = '  +3.17%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  +0.01%  '

# Row 5
$ws.Cells.Item(5, 4).Value = '''580.29'
$ws.Cells.Item(5, 5).Value = '  +1.84%  '

# Row 6
$ws.Cells.Item(6, 4).Value = '''148.19'
$ws.Cells.Item(6, 5).Value = '  +9.39%  '

# Row 7
$ws.Cells.Item(7, 4).Value = '3.458.25'
$ws.Cells.Item(7, 5).Value = '  +3.18%  '

# Row 8
$ws.Cells.Item(8, 5).Value = '  +0.02%  '

# Row 9
$ws.Cells.Item(9, 4).Value = '''0.473'
$ws.Cells.Item(9, 5).Value = '  +1.35%  '

# Row 10
$ws.Cells.Item(10, 4).Value = '''7.71'
$ws.Cells.Item(10, 5).Value = '  +3.36%  '

# Row 11
$ws.Cells.Item(11, 5).Value = '  +1.96%  '

# Row 12
$ws.Cells.Item(12, 4).Value = '''0.391'
$ws.Cells.Item(12, 5).Value = '  +1.59%  '

# Row 13
$ws.Cells.Item(13, 4).Value = '4.049.60'
$ws.Cells.Item(13, 5).Value = '  +3.25%  '

# Row 14
$ws.Cells.Item(14, 4).Value = '''28.08'
$ws.Cells.Item(14, 5).Value = '  +8.69%  '

# Row 15
$ws.Cells.Item(15, 5).Value = '  -0.12%  '

# Row 16
$ws.Cells.Item(16, 4).Value = '''0.0000175'
$ws.Cells.Item(16, 5).Value = '  +2.13%  '

# Row 17
$ws.Cells.Item(17, 4).Value = '3.461.65'
$ws.Cells.Item(17, 5).Value = '  +3.48%  '

# Row 18
$ws.Cells.Item(18, 4).Value = '61.856.83'
$ws.Cells.Item(18, 5).Value = '  +1.12%  '

# Row 19
$ws.Cells.Item(19, 5).Value = '  +9.51%  '

# Row 20
$ws.Cells.Item(20, 4).Value = '''14.43'
$ws.Cells.Item(20, 5).Value = '  +3.30%  '

# Row 21
$ws.Cells.Item(21, 4).Value = '''9.46'
$ws.Cells.Item(21, 5).Value = '  +2.68%  '

# Row 22
$ws.Cells.Item(22, 4).Value = '''385.53'
$ws.Cells.Item(22, 5).Value = '  +2.24%  '

# Row 23
$ws.Cells.Item(23, 4).Value = '''0.570'
$ws.Cells.Item(23, 5).Value = '  +3.50%  '

# Row 24
$ws.Cells.Item(24, 4).Value = '3.593.83'
$ws.Cells.Item(24, 5).Value = '  +3.12%  '

# Row 25
$ws.Cells.Item(25, 5).Value = '  +0.93%  '

# Row 26
$ws.Cells.Item(26, 5).Value = '  +0.12%  '

# Row 27
$ws.Cells.Item(27, 4).Value = '''72.39'
$ws.Cells.Item(27, 5).Value = '  +2.09%  '

# Row 28
$ws.Cells.Item(28, 5).Value = '  -0.21%  '

# Row 29
$ws.Cells.Item(29, 5).Value = '  +9.41%  '

# Row 30
$ws.Cells.Item(30, 5).Value = '  +5.17%  '

# Row 31
$ws.Cells.Item(31, 2).Value = 'Fetch.AI'
$ws.Cells.Item(31, 3).Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Cells.Item(31, 4).Value = '''1.55'
$ws.Cells.Item(31, 5).Value = '  -12.03%  '

# Row 32
$ws.Cells.Item(32, 2).Value = 'Binance-PegBSC-USD'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Cells.Item(32, 4).Value = '''0.999'
$ws.Cells.Item(32, 5).Value = '  -0.31%  '

# Row 33
$ws.Cells.Item(33, 5).Value = '  +1.84%  '

# Row 34
$ws.Cells.Item(34, 5).Value = '  +2.51%  '

# Row 35
$ws.Cells.Item(35, 5).Value = '  +0.01%  '

# Row 36
$ws.Cells.Item(36, 4).Value = '''24.02'
$ws.Cells.Item(36, 5).Value = '  +2.06%  '

# Row 37
$ws.Cells.Item(37, 5).Value = '  +4.76%  '

# Row 38
$ws.Cells.Item(38, 4).Value = '''5.21'
$ws.Cells.Item(38, 5).Value = '  +0.45%  '

# Row 39
$ws.Cells.Item(39, 5).Value = '  +2.51%  '

# Row 40
$ws.Cells.Item(40, 4).Value = '''166.63'

# Row 41
$ws.Cells.Item(41, 4).Value = '''0.0788'
$ws.Cells.Item(41, 5).Value = '  +4.15%  '

# Row 42
$ws.Cells.Item(42, 4).Value = '''26.10'
$ws.Cells.Item(42, 5).Value = '  +10.61%  '

# Row 43
$ws.Cells.Item(43, 4).Value = '''0.797'
$ws.Cells.Item(43, 5).Value = '  +3.98%  '

# Row 44
$ws.Cells.Item(44, 2).Value = 'FirstDigitalUSD'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Cells.Item(44, 4).Value = '''1.00'
$ws.Cells.Item(44, 5).Value = '  +0.04%  '

# Row 45
$ws.Cells.Item(45, 2).Value = 'Stacks'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Cells.Item(45, 4).Value = '''1.73'
$ws.Cells.Item(45, 5).Value = '  +1.62%  '

# Row 46
$ws.Cells.Item(46, 5).Value = '  +2.77%  '

# Row 47
$ws.Cells.Item(47, 4).Value = '''42.32'
$ws.Cells.Item(47, 5).Value = '  +2.24%  '

# Row 48
$ws.Cells.Item(48, 2).Value = 'ONDO'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Cells.Item(48, 4).Value = '''1.17'
$ws.Cells.Item(48, 5).Value = '  -2.06%  '

# Row 49
$ws.Cells.Item(49, 2).Value = 'Maker'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Cells.Item(49, 4).Value = '2.608.51'
$ws.Cells.Item(49, 5).Value = '  +11.31%  '

# Row 50
$ws.Cells.Item(50, 5).Value = '  +2.85%  '

# Row 51
$ws.Cells.Item(51, 4).Value = '''23.46'
$ws.Cells.Item(51, 5).Value = '  +2.15%  '
